$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (396) down to the new rows (397-405)
$ws.Range("A396:O396").Copy() | Out-Null
$ws.Range("A397:O405").PasteSpecial(-4122) | Out-Null

# Row 397
$ws.Range("A397").Value = 395
$ws.Range("B397").Value = 44246.55214344907
$ws.Range("C397").Value = "Turku"
$ws.Range("D397").Value = "31-35 v"
$ws.Range("E397").Value = "mies"
$ws.Range("F397").Value = 8
$ws.Range("G397").Value = "Työntekijä / palkollinen"
$ws.Range("H397").Value = 1
$ws.Range("I397").Value = "Senior Software Engineer (Backend)"
$ws.Range("J397").Value = "Etä"
$ws.Range("K397").Value = 5700
$ws.Range("L397").Value = 74100
$ws.Range("M397").Value = $true
$ws.Range("O397").Value = "Ennen koronaa oli osittainen etätyö, koronan jälkeen 100%"

# Row 398
$ws.Range("A398").Value = 396
$ws.Range("B398").Value = 44246.55232758102
$ws.Range("C398").Value = "EU"
$ws.Range("D398").Value = "31-35 v"
$ws.Range("E398").Value = "mies"
$ws.Range("F398").Value = 3
$ws.Range("G398").Value = "Työntekijä / palkollinen"
$ws.Range("H398").Value = 1
$ws.Range("I398").Value = "Ohjelmistokehittäjä"
$ws.Range("K398").Value = 3200
$ws.Range("L398").Value = 40000
$ws.Range("M398").Value = $false

# Row 399
$ws.Range("A399").Value = 397
$ws.Range("B399").Value = 44246.56990268519
$ws.Range("C399").Value = "Jyväskylä"
$ws.Range("D399").Value = "31-35 v"
$ws.Range("E399").Value = "mies"
$ws.Range("F399").Value = 6
$ws.Range("G399").Value = "Työntekijä / palkollinen"
$ws.Range("H399").Value = 1
$ws.Range("I399").Value = "WordPress / Frontend-koodari"
$ws.Range("J399").Value = "Etä"
$ws.Range("K399").Value = 3000
$ws.Range("L399").Value = 37500
$ws.Range("M399").Value = $true

# Row 400
$ws.Range("A400").Value = 398
$ws.Range("B400").Value = 44246.58209564815
$ws.Range("C400").Value = "Jyväskylä"
$ws.Range("D400").Value = "21-25 v"
$ws.Range("E400").Value = "mies"
$ws.Range("F400").Value = 21
$ws.Range("G400").Value = "Työntekijä / palkollinen"
$ws.Range("H400").Value = 1
$ws.Range("I400").Value = "Arkkitehti"
$ws.Range("J400").Value = "50/50"
$ws.Range("K400").Value = 5500
$ws.Range("L400").Value = 75000
$ws.Range("M400").Value = $true

# Row 401
$ws.Range("A401").Value = 399
$ws.Range("B401").Value = 44246.58394819444
$ws.Range("C401").Value = "PK-Seutu"
$ws.Range("D401").Value = "31-35 v"
$ws.Range("E401").Value = "mies"
$ws.Range("F401").Value = 7
$ws.Range("G401").Value = "Työntekijä / palkollinen"
$ws.Range("H401").Value = 1
$ws.Range("I401").Value = "Ohjelmistokehittäjä, backend"
$ws.Range("J401").Value = "50/50"
$ws.Range("K401").Value = 5470
$ws.Range("L401").Value = 94000
$ws.Range("M401").Value = $true

# Row 402
$ws.Range("A402").Value = 400
$ws.Range("B402").Value = 44246.58876788194
$ws.Range("C402").Value = "PK-Seutu"
$ws.Range("D402").Value = "31-35 v"
$ws.Range("F402").Value = 3
$ws.Range("G402").Value = "Työntekijä / palkollinen"
$ws.Range("H402").Value = 1
$ws.Range("I402").Value = "Full stack developer"
$ws.Range("J402").Value = "Etä"
$ws.Range("K402").Value = 5300
$ws.Range("L402").Value = 66250
$ws.Range("M402").Value = $true

# Row 403
$ws.Range("A403").Value = 401
$ws.Range("B403").Value = 44246.59057026621
$ws.Range("C403").Value = "Kuopio"
$ws.Range("D403").Value = "31-35 v"
$ws.Range("E403").Value = "mies"
$ws.Range("F403").Value = 9
$ws.Range("G403").Value = "Työntekijä / palkollinen"
$ws.Range("H403").Value = 0.8
$ws.Range("I403").Value = "CTO"
$ws.Range("J403").Value = "Etä"
$ws.Range("K403").Value = 5200
$ws.Range("L403").Value = 65000
$ws.Range("M403").Value = $true

# Row 404
$ws.Range("A404").Value = 402
$ws.Range("B404").Value = 44246.59166526621
$ws.Range("C404").Value = "PK-Seutu"
$ws.Range("D404").Value = "36-40 v"
$ws.Range("E404").Value = "mies"
$ws.Range("F404").Value = 14
$ws.Range("G404").Value = "Työntekijä / palkollinen"
$ws.Range("H404").Value = 1
$ws.Range("I404").Value = "Projektipäällikkö"
$ws.Range("J404").Value = "50/50"
$ws.Range("K404").Value = 6400
$ws.Range("L404").Value = 102000
$ws.Range("M404").Value = $true

# Row 405
$ws.Range("A405").Value = 403
$ws.Range("B405").Value = 44246.60321177083
$ws.Range("C405").Value = "PK-Seutu"
$ws.Range("D405").Value = "36-40 v"
$ws.Range("E405").Value = "mies"
$ws.Range("F405").Value = 15
$ws.Range("G405").Value = "Työntekijä / palkollinen"
$ws.Range("H405").Value = 1
$ws.Range("I405").Value = "Frontend & UX"
$ws.Range("K405").Value = 5000
$ws.Range("L405").Value = "Optiot"
$ws.Range("M405").Value = $false
